$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap the header values of B1 and C1 ("Apellido" <-> "Nombre")
$b1 = $ws.Range("B1").Value2
$c1 = $ws.Range("C1").Value2
$ws.Range("B1").Value = $c1
$ws.Range("C1").Value = $b1

# Update the selected cell to C9
$ws.Range("C9").Select()
